$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Number of Videos for Train and Test - P4 (F/M) and P5 (G/N) columns filled in
$data = @(
    @{Row=2; F=25; G=24; M=6; N=6}
    @{Row=3; F=28; G=24; M=6; N=6}
    @{Row=4; F=36; G=24; M=12; N=6}
    @{Row=5; F=25; G=24; M=6; N=6}
    @{Row=6; F=23; G=24; M=8; N=6}
    @{Row=7; F=26; G=24; M=6; N=6}
    @{Row=8; F=25; G=24; M=6; N=6}
    @{Row=9; F=28; G=24; M=6; N=6}
    @{Row=10; F=24; G=24; M=6; N=6}
    @{Row=11; F=27; G=24; M=6; N=6}
    @{Row=12; F=25; G=24; M=6; N=6}
    @{Row=13; F=25; G=24; M=6; N=6}
    @{Row=14; F=25; G=24; M=6; N=6}
    @{Row=15; F=24; G=24; M=6; N=6}
    @{Row=16; F=25; G=24; M=6; N=6}
    @{Row=17; F=25; G=24; M=6; N=6}
    @{Row=18; F=27; G=24; M=6; N=6}
    @{Row=19; F=24; G=24; M=6; N=6}
    @{Row=20; F=25; G=24; M=6; N=6}
    @{Row=21; F=25; G=24; M=6; N=6}
    @{Row=22; F=24; G=24; M=6; N=6}
    @{Row=23; F=24; G=24; M=6; N=6}
    @{Row=24; F=24; G=24; M=6; N=6}
    @{Row=25; F=24; G=24; M=6; N=6}
    @{Row=26; F=25; G=24; M=6; N=6}
    @{Row=27; F=25; G=24; M=6; N=6}
    @{Row=28; F=24; G=24; M=6; N=6}
    @{Row=29; F=24; G=24; M=6; N=6}
    @{Row=30; F=24; G=24; M=6; N=6}
    @{Row=31; F=24; G=24; M=6; N=6}
    @{Row=32; F=25; G=24; M=6; N=6}
    @{Row=33; F=24; G=24; M=6; N=6}
    @{Row=34; F=24; G=24; M=6; N=6}
    @{Row=35; F=24; G=24; M=6; N=6}
    @{Row=36; F=25; G=24; M=6; N=6}
    @{Row=37; F=24; G=24; M=6; N=6}
    @{Row=38; F=24; G=24; M=6; N=6}
    @{Row=39; F=24; G=24; M=6; N=6}
    @{Row=40; F=24; G=24; M=6; N=6}
    @{Row=41; F=24; G=24; M=6; N=6}
    @{Row=42; F=24; G=24; M=6; N=6}
    @{Row=43; F=24; G=23; M=6; N=6}
    @{Row=44; F=25; G=24; M=6; N=6}
    @{Row=45; F=24; G=24; M=6; N=6}
    @{Row=46; F=24; G=24; M=6; N=6}
    @{Row=47; F=24; G=24; M=6; N=6}
    @{Row=48; F=24; G=24; M=6; N=6}
    @{Row=49; F=24; G=24; M=6; N=6}
    @{Row=50; F=24; G=24; M=6; N=6}
    @{Row=51; F=24; G=24; M=6; N=6}
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Range("F$r").Value = $entry.F
    $ws.Range("G$r").Value = $entry.G
    $ws.Range("M$r").Value = $entry.M
    $ws.Range("N$r").Value = $entry.N
}

# Update the active selection to match the post-edit workbook state
$ws.Range("M52").Select()
